$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) sometimes contains values that look numeric
# (e.g. "543.73" or "59.380.10"). Excel would otherwise coerce these into
# actual numbers, losing the original text formatting/precision. Temporarily
# mark them as Text, assign the value, then restore the original cell style
# so no unintended formatting changes are introduced.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D29", "D30", "D31", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
$origStyles = @{}
foreach ($cell in $priceCells) {
    $origStyles[$cell] = $ws.Range($cell).Style
    $ws.Range($cell).NumberFormat = "@"
}

# Apply updated cell values (prices, volumes, coin names/links as changed)
$ws.Range("D2").Value = "59.380.10"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.524.35"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "543.73"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "145.20"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "2.547.10"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "5.61"
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("D13").Value = "0.362"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "2.968.62"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "23.74"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "59.283.02"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "2.534.74"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "11.20"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "325.94"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").Value = "5.85"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").Value = "62.11"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "0.434"
$ws.Range("E25").Value = "  -3.46%  "
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("D29").Value = "0.0₃0787"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").Value = "1.83"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").Value = "6.73"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("E32").Value = "  -5.84%  "
$ws.Range("E33").Value = "  +5.75%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "158.77"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "18.76"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "4.39"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  -5.81%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "5.63"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "36.96"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("D41").Value = "3.71"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "294.96"
$ws.Range("E42").Value = "  -6.16%  "
$ws.Range("D43").Value = "0.827"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "0.602"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").Value = "10.82"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "0.0936"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").Value = "18.79"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").Value = "122.84"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").Value = "0.0516"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  -0.72%  "

# Restore original styles on the price cells
foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = $origStyles[$cell]
}

Write-Output "Updated cryptos list"
